# Fixed flow of new student
# Appends 19 new student records (StudentName, Father Name, Mobile Number,
# Alternate Mobile Number, Email, RRNO) to the StudentData sheet, starting
# right after the last existing row (44), matching the sheet's existing
# all-text layout (so phone-number-looking values stay text, not numbers).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$students = @(
    @('Kessler Runte', 'Olin', '6414338856', '9393812098', 'KesslerRunte@yopmail.com', 'AD2646'),
    @('Howell Koelpin', 'Ed', '9575201245', '8706760528', 'HowellKoelpin@yopmail.com', 'AD9829'),
    @('Hills Watsica', 'Rupert', '6449026896', '5378548338', 'HillsWatsica@yopmail.com', 'AD9378'),
    @('Botsford Littel', 'Winford', '8486967261', '5503420437', 'BotsfordLittel@yopmail.com', 'AD2180'),
    @('Turcotte Kuhlman', 'Robert', '8668617344', '7609915767', 'TurcotteKuhlman@yopmail.com', 'AD9518'),
    @('Altenwerth Sipes', 'Jeannine', '9732936558', '6084743216', 'AltenwerthSipes@yopmail.com', 'AD1857'),
    @('Nienow Green', 'Rodolfo', '3194003211', '3004825105', 'NienowGreen@yopmail.com', 'AD3722'),
    @('Rolfson Wehner', 'Lyndon', '5206075466', '3670515960', 'RolfsonWehner@yopmail.com', 'AD3331'),
    @('Wilkinson Koepp', 'Elijah', '5934780039', '7728210581', 'WilkinsonKoepp@yopmail.com', 'AD9912'),
    @('Cruickshank D''Amore', 'Delena', '5212992532', '6081417927', 'CruickshankD''Amore@yopmail.com', 'AD6931'),
    @('Rice Kub', 'Carli', '5726536221', '4596780178', 'RiceKub@yopmail.com', 'AD8647'),
    @('Nader Spencer', 'Richard', '3453608322', '5671299484', 'NaderSpencer@yopmail.com', 'AD1108'),
    @('Flatley McCullough', 'Matthew', '7340007559', '4176087803', 'FlatleyMcCullough@yopmail.com', 'AD6923'),
    @('Renner Littel', 'Christinia', '5697404055', '3623957549', 'RennerLittel@yopmail.com', 'AD5805'),
    @('Kerluke Tromp', 'Keith', '3067104548', '7530749469', 'KerlukeTromp@yopmail.com', 'AD1604'),
    @('Becker Rodriguez', 'Eloise', '4374465510', '3978001471', 'BeckerRodriguez@yopmail.com', 'AD9714'),
    @('Toy Dickens', 'Silas', '4524178908', '7130580017', 'ToyDickens@yopmail.com', 'AD3475'),
    @('O''Conner Fisher', 'Dorie', '5347643850', '6211310336', 'O''ConnerFisher@yopmail.com', 'AD5529'),
    @('Connelly Farrell', 'Theodore', '3631752501', '6181885702', 'ConnellyFarrell@yopmail.com', 'AD1697'),
)

$startRow = 45
for ($i = 0; $i -lt $students.Count; $i++) {
    $row = $startRow + $i
    $values = $students[$i]
    for ($c = 0; $c -lt $values.Count; $c++) {
        $cell = $ws.Cells.Item($row, $c + 1)
        # Force text storage (columns C/D are numeric-looking mobile
        # numbers that must be kept as text, matching the rest of the
        # sheet) without leaving the cell's number format changed.
        $cell.NumberFormat = '@'
        $cell.Value = $values[$c]
        $cell.ClearFormats()
    }
}
